{"js": "// Replace the outdated \"Constelaci\u00f3n de Tauro\" campaign-dates sentence with\n// the updated wording, everywhere it appears in the document body.\nconst oldText = \"Datas da campa\u00f1a de Constelaci\u00f3n de Tauro 2022: 16-25 de xaneiro\";\nconst newText = \"Datas da campa\u00f1a de 2022 que usan Constelaci\u00f3n de Tauro: 16-25 de xaneiro\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the outdated \"Constelaci\u00f3n de Tauro\" campaign-dates sentence with\n# the updated wording, everywhere it appears in the document body.\n$d = $word.ActiveDocument\n\n$oldText = \"Datas da campa\u00f1a de Constelaci\u00f3n de Tauro 2022: 16-25 de xaneiro\"\n$newText = \"Datas da campa\u00f1a de 2022 que usan Constelaci\u00f3n de Tauro: 16-25 de xaneiro\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
